$d = $word.ActiveDocument

# --- Title: "La Palma Earthquakes" -> "Manuscript 1" ---
$d.Content.Find.Execute("La Palma Earthquakes", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Manuscript 1", 2)

# --- Author 1: "Steve Purves" -> "Andreas Ludvig Ohm Svendsen" ---
$d.Content.Find.Execute("Steve Purves", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Andreas Ludvig Ohm Svendsen", 2)

# --- Author 2: "Rowan Cockett" -> "Tore B. Stage" ---
$d.Content.Find.Execute("Rowan Cockett", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tore B. Stage", 2)

# --- Abstract body text -> "This is an abstract" ---
$abstractOld = "In September 2021, a significant jump in seismic activity on the island of La Palma (Canary Islands, Spain) signaled the start of a volcanic crisis that still continues at the time of writing. Earthquake data is continually collected and published by the Instituto Geogr" + [char]0x00E1 + "phico Nacional (IGN)."
$d.Content.Find.Execute($abstractOld, $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is an abstract", 2)
